# Login.xlsx - data
#
# LoginPage: modified locator for login button (commit message) goes with
# refreshed stored test credentials in this workbook: the username (A2)
# and password (B2) cells are updated, the stale "display" override on
# the A2 mailto hyperlink is dropped (kept pointing at the same address),
# and the sheet's saved active-cell selection moves to D7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the stored username (A2) and password (B2) values.
$ws.Range("A2").Value = "balajee.cs@gmail.com"
$ws.Range("B2").Value = "testleafsf1"

# 2. Recreate the A2 hyperlink so it no longer carries an explicit
#    "display text" override, while still pointing at the same mailto
#    address. Re-apply the Hyperlink cell style afterwards so A2's
#    formatting matches what it had before.
$ws.Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:hari.radhakrishnan@qeagle.com") | Out-Null
$ws.Range("A2").Style = "Hyperlink"

# 3. Move the active selection cell to D7, matching the saved sheet view.
$ws.Range("D7").Select() | Out-Null
